$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3; all existing rows 3..29 shift down to 4..30,
# and formatting (e.g. the date style on column D) is carried along by Excel.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new weekly record.
$ws.Range("A3").Value = 6
$ws.Range("B3").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C3").Value = "Metropolitana"
$ws.Range("D3").Value = 44616
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100102
$ws.Range("H3").Value = "Cítricos"
$ws.Range("I3").Value = 100102006
$ws.Range("J3").Value = "Pomelo"
$ws.Range("K3").Value = "Start Ruby"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 24
$ws.Range("N3").Value = 200000
$ws.Range("O3").Value = 200000
$ws.Range("P3").Value = 200000
$ws.Range("Q3").Value = "$/bins (350 kilos)"
$ws.Range("R3").Value = "Región Metropolitana"
$ws.Range("S3").Value = 571
$ws.Range("T3").Value = 350
